$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 45
$ws.Range("A45").Value = "Record"
$ws.Range("B45").Value = "Balanço Geral"
$ws.Range("C45").Value = "Defesa Civil"
$ws.Range("D45").Value = "2025-04-04T12:02"
$ws.Range("E45").Value = "Neutro"
$ws.Range("F45").Value = "Chuva chega a Campos mas ainda não preocupa autoridades. Repórter *ao vivo*. Imagens de chuva. Previsão para sábado é de chuva forte. Registros de precipitação intensa no Porto do Açu. "

# Row 46
$ws.Range("A46").Value = "Record"
$ws.Range("B46").Value = "Balanço Geral"
$ws.Range("C46").Value = "Saúde"
$ws.Range("D46").Value = "2025-04-04T13:14"
$ws.Range("E46").Value = "Positivo"
$ws.Range("F46").Value = "Vacinação contra o vírus influenza começa na próxima segunda-feira. Repórter *ao vivo*. Campanha será realizada em uma única etapa. De acordo com Secretaria de Saúde, vão ser divulgados os locais até o fim desta semana. Crianças de 5 anos, gestantes, idosos e, este ano, funcionários dos Correios e da área administrativa fazem parte do público alvo. "
